# "Generate Report for Handoff" - refreshes the localization-status report
# after a handoff: status text flips from "In Translation" to
# "Ready for handoff" and the associated timestamps advance. Excel then
# widens the (now longer) Status-ish columns to fit the new text.

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-09-05 17:08:11"

# --- zh-cn sheet ------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-09-05 17:08:00"

# --- de-de sheet ------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-09-05 17:08:11"

# --- Widen the columns that now hold the longer "Ready for handoff"
# status text (mirrors Excel auto-resizing the Status columns after the
# text changed).
$wsOverview.Columns.Item(5).ColumnWidth = 16.84
$wsOverview.Columns.Item(6).ColumnWidth = 16.84
$wsZhCn.Columns.Item(3).ColumnWidth = 16.84
$wsDeDe.Columns.Item(3).ColumnWidth = 16.84
